# Case_6_33 diagnostic.xlsx edit:
#   B1 = 0, A2 = 0, B2 = "disconnected_elements"
#   B1 and A2 get a bold, centered (horizontal center / vertical top), thin-boxed style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values -----------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- style --------------------------------------------------------------
# Build the full style (bold font, thin box border, centered/top aligned) on
# B1 first, then clone it onto A2 via Copy/PasteSpecial(xlPasteFormats) so
# both cells end up sharing exactly one cellXf (avoids creating duplicate /
# stray style entries when setting the same properties cell-by-cell).
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Borders.LineStyle = 1   # xlContinuous
$ws.Range("B1").Borders.Weight = 2      # xlThin
$ws.Range("B1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B1").VerticalAlignment = -4160    # xlTop

$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)     # xlPasteFormats
$excel.CutCopyMode = $false
